# Update gh-pages output data (values refreshed at commit 456a3b4)
$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 331
    $ws.Range("F4").Value = 63
}
